# lagledare_truck.xlsx edit
# - Rename "Logistik" category label to "Matförsörj." for the Materialförsörjning rows (A2:A6)
# - Rename "Administrativt" category label to "Administ." for the Administrativt rows (A7:A10)
# - Unhide column A
# - Update sheet selection / scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unhide column A (was hidden="1")
$ws.Columns("A").Hidden = $false

# Update the short category labels in column A.
$ws.Range("A2:A6").Value = "Matförsörj."
$ws.Range("A7:A10").Value = "Administ."

# Update selection / scroll so the saved view matches (was topLeftCell B1 / selection B27)
$ws.Range("E12").Select()
